$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Insert a new (blank) column before column N, shifting the old "Late" (N)
# and "Outstanding" (O -> previously unused / P) columns one place to the
# right, turning the sheet's 16-column layout (A:P) into a 17-column one
# (A:Q) with a new, currently-empty column N.
$ws.Columns("N").Insert() | Out-Null

# Match the author's explicit width for the newly inserted column N.
$ws.Columns("N").ColumnWidth = 9.1

# Restore the active selection recorded for this sheet after the edit.
$ws.Range("U8").Select() | Out-Null
